# TC01_Canine_Filter_FileAssoc-diagnosis.xlsx - icdc regression suite; INS project
#
# 1) Update the Cypher query in B2 (startup sheet) to also project
#    demo.weight AS weight alongside the existing age projection.
# 2) Scroll the sheet so row 2 becomes the first visible row (topLeftCell A2).
# 3) The row got much taller once the query text grew, so size row 2 to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 2)
$query = $cell.Value2
$oldWith = "WITH DISTINCT c, s, demo, diag, co,demo.patient_age_at_enrollment AS age"
$newWith = "WITH DISTINCT c, s, demo, diag, co,demo.patient_age_at_enrollment AS age, demo.weight as weight"
$cell.Value2 = $query.Replace($oldWith, $newWith)

# Scroll so row 2 is the top-left visible row of the window.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1

# The B2 cell wraps text; the longer query text needs a taller row.
$ws.Rows.Item(2).RowHeight = 375
